# TAK update, Psychiatry addition
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Largest" / TAK date update (I4/J4 block) ---
$ws.Range("J4").Value = 45054

# --- Insert "Gastroenterology" ahead of "Genetic Diseases" in the F column list ---
# Push the existing "Genetic Diseases" entry down into F6 (previously empty)...
$ws.Range("F6").Value = $ws.Range("F5").Value()
# ...and replace F5 with the newly added disease category.
$ws.Range("F5").Value = "Gastroenterology"

# --- Psychiatry (F12) update date added in G12, matching the date style used in column G ---
$ws.Range("G11").Copy()
$ws.Range("G12").PasteSpecial(-4122)
$ws.Range("G12").Value = 45067

# --- Update the active selection to reflect where the user last clicked ---
$ws.Range("G13").Select() | Out-Null
